# Updates the cryptos list worksheet with refreshed prices / 1h volume
# percentages (and a coin-row swap / replacement), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# D-column "Price" cells are plain text in the source data (they use
# "."-grouped formatting, e.g. "66.944.77", and values like "1.00" must
# keep their trailing zero), so each Price write is wrapped with a
# text NumberFormat to stop Excel from re-interpreting it as a number;
# ClearFormats() afterwards restores the cell's original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.944.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.597.50"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.596.50"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("E10").Value = "  -4.93%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.42"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.071.99"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("E16").Value = "  -5.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.828.98"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.596.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "363.87"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("E21").Value = "  -5.25%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "67.55"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.727.56"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "581.58"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -4.40%  "
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.68"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("E33").Value = "  -3.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  -8.79%  "
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.18"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.81"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.51"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.613"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.549"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.41%  "
